## Edit: swap the deck's theme palette from "Integral" to the stock
## "Office Theme" colours, and switch the Component 3 sources-of-finance
## table (slide 6) onto a different built-in table style.

$p = $ppt.ActivePresentation

# --- 1. Table style: slide 6, the graphicFrame holding the table ---
$tableSlide = $p.Slides.Item(6)
$tableShape = $tableSlide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{2C5F4924-A5AA-4B5A-93C6-EB5FF29F6AFA}")

# --- 2. Theme colours: Integral -> Office Theme, in clrScheme order ---
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) as 0x00BBGGRR ints.
$master = $p.Slides.Item(1).Master
$themeColors = $master.Theme.ThemeColorScheme

$themeColors.Colors(1).RGB  = 0        # dk1     000000
$themeColors.Colors(2).RGB  = 16777215 # lt1     FFFFFF
$themeColors.Colors(3).RGB  = 6968388  # dk2     44546A
$themeColors.Colors(4).RGB  = 15132391 # lt2     E7E6E6
$themeColors.Colors(5).RGB  = 13998939 # accent1 5B9BD5
$themeColors.Colors(6).RGB  = 3243501  # accent2 ED7D31
$themeColors.Colors(7).RGB  = 10855845 # accent3 A5A5A5
$themeColors.Colors(8).RGB  = 49407    # accent4 FFC000
$themeColors.Colors(9).RGB  = 12874308 # accent5 4472C4
$themeColors.Colors(10).RGB = 4697456  # accent6 70AD47
$themeColors.Colors(11).RGB = 12673797 # hlink   0563C1
$themeColors.Colors(12).RGB = 7491477  # folHlink 954F72
